$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder email lists in column G (Recorded By) for the rows listed in the commit diff.
# Using Cells.Item(row, 7) since column G is the 7th column.

$ws.Cells.Item(3, 7).Value = 'servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Cells.Item(4, 7).Value = 'Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Cells.Item(9, 7).Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Cells.Item(10, 7).Value = 'aml.awwad@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Cells.Item(15, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(16, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(27, 7).Value = 'marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, System'
$ws.Cells.Item(29, 7).Value = 'marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(34, 7).Value = 'servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, shaimaa.ahmed@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Cells.Item(35, 7).Value = 'asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg'
$ws.Cells.Item(38, 7).Value = 'AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg'
$ws.Cells.Item(40, 7).Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg'
$ws.Cells.Item(41, 7).Value = 'aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(46, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Cells.Item(58, 7).Value = 'marian.samir@med.asu.edu.eg, randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, System'
$ws.Cells.Item(60, 7).Value = 'marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(64, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg'
$ws.Cells.Item(65, 7).Value = 'majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(66, 7).Value = 'Veronia.rafat@med.asu.edu.eg, hananragab@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg'
$ws.Cells.Item(71, 7).Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Cells.Item(72, 7).Value = 'aml.awwad@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, aya.saeed@med.asu.edu.eg'
$ws.Cells.Item(76, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(77, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Cells.Item(89, 7).Value = 'nourhan.mostafa@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg'
$ws.Cells.Item(91, 7).Value = 'marina_atef@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, nahla91@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg'
$ws.Cells.Item(95, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Cells.Item(96, 7).Value = 'servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(97, 7).Value = 'asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(102, 7).Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Cells.Item(103, 7).Value = 'aml.awwad@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg'
$ws.Cells.Item(108, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(109, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg'
$ws.Cells.Item(120, 7).Value = 'marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg'
$ws.Cells.Item(122, 7).Value = 'Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(126, 7).Value = 'mennatulla.medhat@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Cells.Item(127, 7).Value = 'servinaz@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm'
$ws.Cells.Item(128, 7).Value = 'asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(133, 7).Value = 'aml.awwad@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Cells.Item(138, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(139, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(140, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Cells.Item(151, 7).Value = 'Amr-Saeed@med.asu.edu.eg, takyas@med.asu.edu.eg, enas.omran@med.asu.edu.eg, user@user.com'
$ws.Cells.Item(153, 7).Value = 'Monica.Eshak@med.asu.edu.eg, eman.samir@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(154, 7).Value = 'marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
$ws.Cells.Item(157, 7).Value = 'majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, servinaz@med.asu.edu.eg'
$ws.Cells.Item(158, 7).Value = 'majorelle.magdy@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg'
$ws.Cells.Item(159, 7).Value = 'asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, gehanadel@med.asu.edu.eg'
$ws.Cells.Item(164, 7).Value = 'Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Cells.Item(170, 7).Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg'
$ws.Cells.Item(182, 7).Value = 'marian.samir@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, dr_ahmed_mostafa@med.asu.edu.eg'
$ws.Cells.Item(184, 7).Value = 'Monica.Eshak@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, nahla91@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Cells.Item(185, 7).Value = 'marina_atef@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg'
